$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values; prefix with an apostrophe so Excel
# keeps these as literal text instead of auto-converting to numbers
# (these values use "." as a thousands-style separator, not a decimal point).
$ws.Range('D2').Value = "'26.346.83"
$ws.Range('D3').Value = "'1.792.56"
$ws.Range('D4').Value = "'1.009"
$ws.Range('D5').Value = "'1.008"
$ws.Range('D6').Value = "'307.88"
$ws.Range('D7').Value = "'0.4527"
$ws.Range('D8').Value = "'0.3588"
$ws.Range('D9').Value = "'45.97"
$ws.Range('D10').Value = "'0.07110"
$ws.Range('D11').Value = "'0.8847"
$ws.Range('D12').Value = "'0.07817"
$ws.Range('D13').Value = "'19.50"
$ws.Range('D14').Value = "'1.786.77"
$ws.Range('D15').Value = "'5.282"
$ws.Range('D16').Value = "'6.323"
$ws.Range('D17').Value = "'84.80"
$ws.Range('D18').Value = "'1.010"
$ws.Range('D19').Value = "'0.000008542"
$ws.Range('D21').Value = "'14.27"
$ws.Range('D22').Value = "'26.383.00"
$ws.Range('D23').Value = "'4.988"
$ws.Range('D24').Value = "'10.50"
$ws.Range('D25').Value = "'2.001.58"
$ws.Range('D26').Value = "'1.967"
$ws.Range('D27').Value = "'152.32"
$ws.Range('D28').Value = "'17.89"
$ws.Range('D29').Value = "'2.033"
$ws.Range('D30').Value = "'111.96"
$ws.Range('D31').Value = "'4.865"
$ws.Range('D32').Value = "'0.08659"
$ws.Range('D33').Value = "'3.047"
$ws.Range('D34').Value = "'4.449"
$ws.Range('D35').Value = "'0.7252"
$ws.Range('D36').Value = "'2.721"
$ws.Range('D37').Value = "'1.108"
$ws.Range('D38').Value = "'1.073"
$ws.Range('D39').Value = "'0.01933"
$ws.Range('D40').Value = "'0.05099"
$ws.Range('D41').Value = "'2.877"
$ws.Range('D42').Value = "'0.5113"
$ws.Range('D43').Value = "'6.871"
$ws.Range('D44').Value = "'0.1515"
$ws.Range('D45').Value = "'7.991"
$ws.Range('D47').Value = "'0.4640"
$ws.Range('D50').Value = "'1.581"
$ws.Range('D51').Value = "'0.05971"

# --- Update Volume(1h) (column E) values. These already contain non-numeric
# characters (%, spaces) so Excel stores them as text without extra handling.
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('E3').Value = '  -2.07%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E7').Value = '  -2.08%  '
$ws.Range('E8').Value = '  -3.08%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('E14').Value = '  -2.81%  '
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('E19').Value = '  -2.15%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('E25').Value = '  -3.23%  '
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('E28').Value = '  -1.88%  '
$ws.Range('E29').Value = '  +3.46%  '
$ws.Range('E30').Value = '  -1.36%  '
$ws.Range('E31').Value = '  -1.31%  '
$ws.Range('E33').Value = '  -2.69%  '
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').Value = '  -3.86%  '
$ws.Range('E36').Value = '  +6.07%  '
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('E38').Value = '  -1.56%  '
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('E40').Value = '  -0.51%  '
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('E42').Value = '  +2.78%  '
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('E44').Value = '  -5.11%  '
$ws.Range('E45').Value = '  -4.04%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('E51').Value = '  -2.10%  '


# --- Rows 48/49: the coin that ranked 47th/48th changed. The row that used
# to list "Quant" now lists "EnergySwap" (and vice versa, with its own
# updated Price/Volume figures), so update B/C/D/E directly for both rows.
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'9.880"
$ws.Range('E48').Value = '  -2.86%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = "'100.81"
$ws.Range('E49').Value = '  -1.46%  '
